# Auto-update gym prices

$wb = $excel.ActiveWorkbook

# Sheet: "4x4 Squat Racks"
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")

# C2: price update. Force text format first so the literal "$2,065.00"
# string is preserved instead of being auto-parsed into a currency number,
# then restore the cell's original (default) style.
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "$2,065.00"
$ws1.Range("C2").Style = "Normal"

# C3: price no longer available.
$ws1.Range("C3").Value = "Not available"

# Sheet: "Squat Stands"
$ws2 = $wb.Worksheets.Item("Squat Stands")

# C2: price update.
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "$1,494.00"
$ws2.Range("C2").Style = "Normal"

# C3: price no longer available.
$ws2.Range("C3").Value = "Price not available"
